$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the "Status" column (column I) contents - removes header + data values
$ws.Range("I1:I4").ClearContents() | Out-Null

# Fix the Sr. No. value in row 4 (was 2, should be 3)
$ws.Range("A4").Value = 3

# Update the selection to column I (now empty) as seen in the target workbook
$ws.Range("I1:I1048576").Select() | Out-Null

$ws.Activate()
